$wb = $excel.ActiveWorkbook

# ---- Summary ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = -173776.4612130868
$ws.Range("B7").Value = 5413711.842050619
$ws.Range("B8").Value = 22726010.95505212
$ws.Range("B10").Value = 4288768.643654695

# ---- Fed-in Capacity ----
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("J44").Value = 169.0966151720738
$ws.Range("K44").Value = 220.0898510449805
$ws.Range("L44").Value = 235.7664149699872
$ws.Range("M44").Value = 230.3462332272727
$ws.Range("N44").Value = 229.4130635965909
$ws.Range("O44").Value = 230.0982114216867
$ws.Range("P44").Value = 231.2329957552695
$ws.Range("Q44").Value = 212.3149906599047
$ws.Range("J45").Value = 126.0910353404088
$ws.Range("K45").Value = 137.841438974359
$ws.Range("L45").Value = 138.5543797798742
$ws.Range("M45").Value = 142.1340339220183
$ws.Range("N45").Value = 131.3417120833333
$ws.Range("O45").Value = 142.5962444444444
$ws.Range("P45").Value = 133.9744074143302
$ws.Range("Q45").Value = 139.9817740860215
$ws.Range("K46").Value = 106.7437663446525
$ws.Range("L46").Value = 134.8846762812383
$ws.Range("M46").Value = 138.9257839476051
$ws.Range("N46").Value = 127.6855444652332
$ws.Range("O46").Value = 138.4565384518428
$ws.Range("P46").Value = 135.0065633140411

# ---- Unmet Demand ----
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("G44").Value = 415.302737515135
$ws.Range("H44").Value = 339.4748021157671
$ws.Range("I44").Value = 210.4758895704059
$ws.Range("J44").Value = 11.94928935461252
$ws.Range("Q44").Value = 9.990699214544804
$ws.Range("R44").Value = 149.8691179411497
$ws.Range("S44").Value = 209.0200695862453
$ws.Range("T44").Value = 223.0958495641314
$ws.Range("U44").Value = 251.3456529078365
$ws.Range("G45").Value = 137.3435171632106
$ws.Range("H45").Value = 112.2354442364965
$ws.Range("I45").Value = 89.39663285141508
$ws.Range("J45").Value = 0.7465913262578567
$ws.Range("R45").Value = 100.1578341526431
$ws.Range("S45").Value = 171.6831711038378
$ws.Range("T45").Value = 200.1647286948216
$ws.Range("U45").Value = 225.9413820809748
$ws.Range("G46").Value = 167.9909793584588
$ws.Range("H46").Value = 162.2271725074396
$ws.Range("I46").Value = 155.4504749272583
$ws.Range("J46").Value = 93.35918011667277
$ws.Range("K46").Value = 22.26949182588285
$ws.Range("P46").Value = 2.721440735106512
$ws.Range("Q46").Value = 86.16204325169439
$ws.Range("R46").Value = 177.2933913771695
$ws.Range("S46").Value = 224.0165980369723
$ws.Range("T46").Value = 227.9455894282815
$ws.Range("U46").Value = 286.3190293564909

# ---- Household Surplus ----
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B16").Value = 370900.54908349

# ---- Costs and Revenues ----
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 43002.96221257856
$ws.Range("D2").Value = 43002.96221257856
$ws.Range("E2").Value = 43002.96221257856
$ws.Range("L2").Value = 43002.96221257856
$ws.Range("P2").Value = 43002.96221257855
$ws.Range("P3").Value = 0
$ws.Range("P4").Value = 37627.59193600624
$ws.Range("P5").Value = 0
$ws.Range("B6").Value = -43402.90369346245
$ws.Range("C6").Value = -43402.90369346245
$ws.Range("D6").Value = -43402.90369346245
$ws.Range("E6").Value = -9775.303693462443
$ws.Range("F6").Value = -9775.303693462451
$ws.Range("G6").Value = -9775.303693462451
$ws.Range("H6").Value = -9775.303693462451
$ws.Range("I6").Value = -9775.303693462451
$ws.Range("J6").Value = -9775.303693462451
$ws.Range("K6").Value = -9775.303693462451
$ws.Range("L6").Value = -9775.303693462443
$ws.Range("M6").Value = -9775.303693462451
$ws.Range("N6").Value = -9775.303693462451
$ws.Range("O6").Value = -9775.303693462451
$ws.Range("P6").Value = -9775.303693462451

# ---- Installed Capacities ----
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Range("P3").Value = 0

# ---- Added Capacities ----
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("P3").Value = 0

# ---- PV Dispatch ----
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G44:U46").Value = 0
